$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header shared-text cells (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# --- Donor cells used as style templates for text placeholder cells ("0" / "***.*") ---
$zeroDonor = $ws.Range("D14")   # style 14, text "0"
$naDonor = $ws.Range("E14")     # style 14, text "***.*"

# --- Cells changing from a number to a text placeholder (copy style+value from donor) ---
$zeroDonor.Copy($ws.Range("C14"))
$zeroDonor.Copy($ws.Range("D15"))
$naDonor.Copy($ws.Range("E15"))
$zeroDonor.Copy($ws.Range("D26"))
$naDonor.Copy($ws.Range("E26"))
$zeroDonor.Copy($ws.Range("C28"))
$zeroDonor.Copy($ws.Range("C29"))

# --- Cells changing from a text placeholder to a number (set NumberFormat then Value) ---
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 2
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("D17").Value = 4
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E17").Value = -25
$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("C20").Value = 4
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("F22").NumberFormat = '#,##0'
$ws.Range("F22").Value = 1
$ws.Range("F30").NumberFormat = '#,##0'
$ws.Range("F30").Value = 1

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F14").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = -11.764705882352
$ws.Range("L16").Value = 7.142857142857
$ws.Range("M16").Value = -52.380952380952
$ws.Range("N16").Value = -91.404011461318
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -18.75
$ws.Range("I17").Value = 54
$ws.Range("J17").Value = 56
$ws.Range("K17").Value = -3.571428571428
$ws.Range("L17").Value = 45.945945945945
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = -61.971830985915
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = -11.111111111111
$ws.Range("L18").Value = 28
$ws.Range("M18").Value = -15.789473684210
$ws.Range("N18").Value = -86.497890295358
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 14.285714285714
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 29.166666666666
$ws.Range("I19").Value = 96
$ws.Range("J19").Value = 89
$ws.Range("K19").Value = 7.865168539325
$ws.Range("L19").Value = 24.675324675324
$ws.Range("M19").Value = -20
$ws.Range("N19").Value = -41.463414634146
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 25
$ws.Range("J20").Value = 24
$ws.Range("K20").Value = 4.166666666666
$ws.Range("L20").Value = 150
$ws.Range("M20").Value = -13.793103448275
$ws.Range("N20").Value = -90.157480314960
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 11.764705882352
$ws.Range("F21").Value = 64
$ws.Range("G21").Value = 64
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 240
$ws.Range("J21").Value = 242
$ws.Range("K21").Value = -0.826446280991
$ws.Range("L21").Value = 34.831460674157
$ws.Range("M21").Value = -17.525773195876
$ws.Range("N21").Value = -79.452054794520
$ws.Range("I22").Value = 4
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -73.333333333333
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("I23").Value = 31
$ws.Range("J23").Value = 30
$ws.Range("K23").Value = 3.333333333333
$ws.Range("L23").Value = 3.333333333333
$ws.Range("M23").Value = 0
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -46.153846153846
$ws.Range("F24").Value = 55
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = -8.333333333333
$ws.Range("I24").Value = 221
$ws.Range("J24").Value = 201
$ws.Range("K24").Value = 9.950248756218
$ws.Range("L24").Value = 40.764331210191
$ws.Range("M24").Value = -25.337837837837
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 80
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 11.111111111111
$ws.Range("I25").Value = 107
$ws.Range("J25").Value = 111
$ws.Range("K25").Value = -3.603603603603
$ws.Range("L25").Value = 48.611111111111
$ws.Range("M25").Value = 16.304347826087
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 128.571428571429
$ws.Range("L27").Value = 45.454545454545
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -66.666666666666
$ws.Range("N28").Value = -84.848484848484
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = -66.666666666666
$ws.Range("N29").Value = -85.714285714285
$ws.Range("I30").Value = 3
$ws.Range("K30").Value = 200
$ws.Range("L30").Value = 50
